# "work on intros and overall data"
#
# - rename the header labels in row 1 to the new snake_case / lowercase
#   column names
# - shrink the header row height (100 -> 40)
# - fill in the previously-blank "Asset Size" values for 1964-1967 with 0
# - move the frozen-pane top-left cell back to A2 and select B2
# - nudge the saved window position (xWindow)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row relabeling -------------------------------------------------
$ws.Range("B1").Value = "asset size"
$ws.Range("C1").Value = "dollars_mn"
$ws.Range("D1").Value = "dollars_nyc"
$ws.Range("E1").Value = "dollars_other"
$ws.Range("F1").Value = "dollars_total"
$ws.Range("G1").Value = "dollars_gp"
$ws.Range("H1").Value = "dollars_nyc_fv"
$ws.Range("I1").Value = "dollars_ts"
$ws.Range("J1").Value = "dollars_mn_fv"
$ws.Range("K1").Value = "applications"
$ws.Range("L1").Value = "approved"

# --- header row height ------------------------------------------------------
$ws.Rows(1).RowHeight = 40

# --- fill in missing Asset Size values for the first 4 data rows -----------
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0

# --- view state: selected cell (frozen pane already anchors at A2) ---------
$null = $ws.Range("B2").Select()

# --- saved window position ---------------------------------------------------
$excel.ActiveWindow.Left = 740
$excel.Left = 740
